$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the two new columns
$ws.Range("P1").Value = "HR_Paired_Pre"
$ws.Range("Q1").Value = "HR_Paired_Post"

# Fill in HR_Paired_Pre / HR_Paired_Post values for each row
$ws.Range("P2").Value = 0.090472618
$ws.Range("Q2").Value = 0.103796438
$ws.Range("P3").Value = 0.24547934799999999
$ws.Range("Q3").Value = 0.214533266
$ws.Range("P4").Value = 0.199352819
$ws.Range("Q4").Value = 0.378937072
$ws.Range("P5").Value = 0.133937669
$ws.Range("Q5").Value = 0.166697307
$ws.Range("P6").Value = 0.127523629
$ws.Range("Q6").Value = 0.366747839
$ws.Range("P7").Value = 0.209374757
$ws.Range("Q7").Value = 0.114302677
$ws.Range("P8").Value = 0.12208941
$ws.Range("Q8").Value = 0.149890458
$ws.Range("P9").Value = 0.143736906
$ws.Range("Q9").Value = 0.267204691
$ws.Range("P10").Value = 0.143348786
$ws.Range("Q10").Value = 0.274748287
$ws.Range("P11").Value = 0.291508123
$ws.Range("Q11").Value = 0.322332322
$ws.Range("P12").Value = 0.156780807
$ws.Range("Q12").Value = 0.484655966
$ws.Range("P13").Value = 0.099652075
$ws.Range("Q13").Value = 0.144480399
$ws.Range("P14").Value = 0.157077178
$ws.Range("Q14").Value = 0.237597541
$ws.Range("P15").Value = 0.901472804
$ws.Range("Q15").Value = 0.84183653
$ws.Range("P16").Value = 0.553186199
$ws.Range("Q16").Value = 0.613760385
$ws.Range("P17").Value = 0.228015143
$ws.Range("Q17").Value = 0.167915071
$ws.Range("P18").Value = 0.118884248
$ws.Range("Q18").Value = 0.167747994
$ws.Range("P19").Value = 0.031441819
$ws.Range("Q19").Value = 0.041805957

# Update the selection to match the saved workbook state
$ws.Range("F2").Select()
